# "break out stock.yaml completed"
# Target sheet is "3 V 0.3" (the chartink screener breakout sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3 V 0.3")

# E2 was stored as text "532900" (bsecode) - it should really be numeric.
$ws.Range("E2").Value = 532900

# Append the new screener row (row 3) with the latest scrape.
$ws.Range("A3").Value = "12/06/2024 06:45:30"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "PAISALO"
$ws.Range("D3").Value = "Paisalo Digital Ltd"
# bsecode keeps coming through as text from the source feed - force text entry
# (leading apostrophe = Excel "quote prefix") so it doesn't become numeric.
$ws.Range("E3").Value = "'532900"
$ws.Range("F3").Value = 9.130000000000001
$ws.Range("G3").Value = 69.36
$ws.Range("H3").Value = 3418345
